$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to make edits, re-protect at the end.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclosure note (A38).
$disclosureText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."
$ws.Range("A38").Value = $disclosureText

# Update Weight (D) and Percent Change (E) values for each holding row.
$ws.Range("D2").Value = 0.03546601036748503  # WFCPRN
$ws.Range("E2").Value = 0
$ws.Range("D3").Value = 0.02028569746926729  # COFPRG
$ws.Range("E3").Value = 0.001174628034455782
$ws.Range("D4").Value = 0.01936312042033669  # SCEPRG
$ws.Range("E4").Value = -0.001598721023181482
$ws.Range("D5").Value = 0.03789999468448985  # MSPRI
$ws.Range("E5").Value = 0.002444987775061191
$ws.Range("D6").Value = 0.03424023539767246  # JPMPRG
$ws.Range("E6").Value = -0.000400000000000067
$ws.Range("D7").Value = 0.01973841376623445  # JPHY
$ws.Range("E7").Value = 0.001741149158444477
$ws.Range("D8").Value = 0.03725719472628076  # IGF
$ws.Range("E8").Value = 0.003600932005930924
$ws.Range("D9").Value = 0.02038039908687275  # EMB
$ws.Range("E9").Value = 0.00107923374404173
$ws.Range("D10").Value = 0.026185506416391  # DUK
$ws.Range("E10").Value = -0.005347073692397397
$ws.Range("D11").Value = 0.02399362220245821  # FDVV
$ws.Range("E11").Value = 0.007692307692307665
$ws.Range("D12").Value = 0.05728623044589799  # FDRR
$ws.Range("E12").Value = 0.006881822496440204
$ws.Range("D13").Value = 0.02479043957639661  # GSPRJ
$ws.Range("E13").Value = 0.002587800369685755
$ws.Range("D14").Value = 0.02687876298914149  # SO
$ws.Range("E14").Value = 0.0006213109661385996
$ws.Range("D15").Value = 0.03245271690794534  # VZ
$ws.Range("E15").Value = 0.001230012300122985
$ws.Range("D16").Value = 0.01973026739052646  # PFE
$ws.Range("E16").Value = -0.003504380475594426
$ws.Range("D17").Value = 0.0312897199459325  # CPT
$ws.Range("E17").Value = 0.01399397933447233
$ws.Range("D18").Value = 0.0419940576262398  # PSK
$ws.Range("E18").Value = 0.00253397834600344
$ws.Range("D19").Value = 0.1259884863198931  # PGX
$ws.Range("E19").Value = -0.001327140013271277
$ws.Range("D20").Value = 0.009272306660537877  # MRK
$ws.Range("E20").Value = -0.001894417782268287
$ws.Range("D21").Value = 0.01543621092514519  # VPU
$ws.Range("E21").Value = -0.0008311976172334923
$ws.Range("D22").Value = 0.01734653602867035  # AVB
$ws.Range("E22").Value = 0.01491941519884254
$ws.Range("D23").Value = 0.01541711785707958  # HTA
$ws.Range("E23").Value = 0.005399568034557101
$ws.Range("D24").Value = 0.02152338651355282  # IYR
$ws.Range("E24").Value = 0.01129081476960647
$ws.Range("D25").Value = 0.0126866054643037  # MLPA
$ws.Range("E25").Value = 0.01075557945684325
$ws.Range("D26").Value = 0.04240723161917974  # VDC
$ws.Range("E26").Value = 0.003727916148900556
$ws.Range("D27").Value = 0.02388456259766741  # MINT
$ws.Range("E27").Value = 0.0002941753284957382
$ws.Range("D28").Value = 0.04552673345201245  # PFXF
$ws.Range("E28").Value = 0.005240590757503716
$ws.Range("D29").Value = 0.05560471667007116  # TDIV
$ws.Range("E29").Value = 0.01238334529791807
$ws.Range("D30").Value = 0.01286373822110445  # APLE
$ws.Range("E30").Value = 0.01511169513797639
$ws.Range("D31").Value = 0.02064495263798993  # VRP
$ws.Range("E31").Value = -0.001529051987767538
$ws.Range("D32").Value = 0.01336514764593163  # MPW
$ws.Range("E32").Value = 0.01476190476190475
$ws.Range("D33").Value = 0.04176020572857212  # PHB
$ws.Range("E33").Value = 0
$ws.Range("D34").Value = 0.01703967223871976  # VGK
$ws.Range("E34").Value = 0.006442166910688218
$ws.Range("D35").Value = 0.9999999999999999  # Total
$ws.Range("E35").Value = 0.003383138916682737

# Restore sheet protection.
$ws.Protect("D382")
